# ---------------------------------------------------------------------------
# C5-PowerPoint.pptx edit
#   1) Slide 6's table switches from the custom "Table_0" style to the
#      built-in PowerPoint table style "Medium Style 2 - Accent 3"
#      ({D563CD09-B2D8-4B5A-8008-60E44CA1D84F}).
#   2) The deck's theme palette (currently the "Integral" theme) is
#      swapped for the stock Office theme colours.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 -------------------------------------------
$s = $p.Slides.Item(6)
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).HasTable) {
        $tableShape = $s.Shapes.Item($i)
        break
    }
}
$tableShape.Table.ApplyStyle("{D563CD09-B2D8-4B5A-8008-60E44CA1D84F}")

# --- 2) Theme colour scheme -------------------------------------------
$theme = $p.Slides.Item(1).Master.Theme
$clr = $theme.ThemeColorScheme

function Set-ThemeRgb($index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $clr.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeRgb 1  "000000"   # dk1
Set-ThemeRgb 2  "FFFFFF"   # lt1
Set-ThemeRgb 3  "44546A"   # dk2
Set-ThemeRgb 4  "E7E6E6"   # lt2
Set-ThemeRgb 5  "5B9BD5"   # accent1
Set-ThemeRgb 6  "ED7D31"   # accent2
Set-ThemeRgb 7  "A5A5A5"   # accent3
Set-ThemeRgb 8  "FFC000"   # accent4
Set-ThemeRgb 9  "4472C4"   # accent5
Set-ThemeRgb 10 "70AD47"   # accent6
Set-ThemeRgb 11 "0563C1"   # hlink
Set-ThemeRgb 12 "954F72"   # folHlink
